# Changed the job id
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rotate the job id values in column A (rows 2-10): the value that was
# in the last row moves to the top, and the rest shift down by one row.
$ws.Range("A2").Value = 9
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8

# Update the active selection to match the saved view state.
$ws.Range("B12").Select()
